# Refresh the "last status check" timestamp shown in the header (F1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Last status check on: 23.02.2022 05:15"

# Row 5 ("Makro") got a fresh price-check reading:
#  - current price (B5) moved 36.7 -> 36.9
#  - previous price is now recorded as the old B5 value (C5) 36.5 -> 36.7
#  - the delta (D5) is now written out as a signed text label instead of a
#    plain number
#  - the "old date" column (E5) is now a plain text timestamp instead of a
#    numeric date serial, so clear its date number-format back to the
#    default/general style once the new text value is in place.
$ws.Range("B5").Value = 36.9
$ws.Range("C5").Value = 36.7

$ws.Range("D5").Value = "'+0.2"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "2022-02-23 05:15:07"
$ws.Range("E5").Style = "Normal"
